# Fruta / hortaliza, semanal
# Update weekly price records for Granada (Mercado Mayorista Lo Valledor de Santiago):
# dates, variedad/calidad, volumen, precios and origen for rows 2-6 and 9-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44721
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 300000
$ws.Range("O2").Value = 300000
$ws.Range("P2").Value = 300000
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 750
# Row 3
$ws.Range("D3").Value = 45043
$ws.Range("K3").Value = "Wonderfull"
$ws.Range("M3").Value = 18
$ws.Range("N3").Value = 300000
$ws.Range("O3").Value = 315000
$ws.Range("P3").Value = 307500
$ws.Range("Q3").Value = "`$/bins (400 kilos)"
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 769
$ws.Range("T3").Value = 400
# Row 4
$ws.Range("D4").Value = 45043
$ws.Range("K4").Value = "Wonderfull"
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 270000
$ws.Range("O4").Value = 270000
$ws.Range("P4").Value = 270000
$ws.Range("Q4").Value = "`$/bins (400 kilos)"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 675
$ws.Range("T4").Value = 400
# Row 5
$ws.Range("D5").Value = 44285
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 8
$ws.Range("N5").Value = 280000
$ws.Range("O5").Value = 300000
$ws.Range("P5").Value = 290000
$ws.Range("Q5").Value = "`$/bins (400 kilos)"
$ws.Range("S5").Value = 725
$ws.Range("T5").Value = 400
# Row 6
$ws.Range("D6").Value = 44280
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 15
$ws.Range("N6").Value = 360000
$ws.Range("O6").Value = 360000
$ws.Range("P6").Value = 360000
$ws.Range("Q6").Value = "`$/bins (450 kilos)"
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 450
# Row 9
$ws.Range("D9").Value = 44662
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range("Q9").Value = "`$/caja 18 kilos granel"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 18
# Row 10
$ws.Range("D10").Value = 44662
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("Q10").Value = "`$/caja 18 kilos granel"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 889
$ws.Range("T10").Value = 18
# Row 11
$ws.Range("D11").Value = 44334
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 16
$ws.Range("N11").Value = 240000
$ws.Range("O11").Value = 250000
$ws.Range("P11").Value = 245000
$ws.Range("Q11").Value = "`$/bins (450 kilos)"
$ws.Range("R11").Value = "Provincia de Limarí"
$ws.Range("S11").Value = 544
$ws.Range("T11").Value = 450
# Row 12
$ws.Range("D12").Value = 44266
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 4800
$ws.Range("O12").Value = 4800
$ws.Range("P12").Value = 4800
$ws.Range("Q12").Value = "`$/bandeja 4 kilos"
$ws.Range("R12").Value = "Provincia del Elquí"
$ws.Range("S12").Value = 1200
$ws.Range("T12").Value = 4
# Row 13
$ws.Range("D13").Value = 44266
$ws.Range("K13").Value = "Wonderfull"
$ws.Range("L13").Value = "Tercera"
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 4000
$ws.Range("O13").Value = 4000
$ws.Range("P13").Value = 4000
$ws.Range("Q13").Value = "`$/bandeja 4 kilos"
$ws.Range("R13").Value = "Provincia del Elquí"
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 4
# Row 14
$ws.Range("D14").Value = 44320
$ws.Range("M14").Value = 12
$ws.Range("N14").Value = 250000
$ws.Range("O14").Value = 260000
$ws.Range("P14").Value = 255000
$ws.Range("Q14").Value = "`$/bins (400 kilos)"
$ws.Range("S14").Value = 638
$ws.Range("T14").Value = 400
# Row 15
$ws.Range("D15").Value = 44307
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17000
$ws.Range("Q15").Value = "`$/caja 15 kilos granel"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1133
$ws.Range("T15").Value = 15
